$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for "Acelga" (Femacal de La
# Calera) ahead of the existing series, so insert a fresh row at row 199 and
# push the rest of the table down by one row.
$ws.Rows("199:199").Insert()

# Populate the newly inserted row with the new observation's data.
$ws.Range("A199").Value = 3
$ws.Range("B199").Value = "Femacal de La Calera"
$ws.Range("C199").Value = "Coquimbo"
$ws.Range("D199").Value = 44522
$ws.Range("E199").Value = 5
$ws.Range("F199").Value = 100112009
$ws.Range("G199").Value = "Acelga"
$ws.Range("H199").Value = "Sin especificar"
$ws.Range("I199").Value = "Primera"
$ws.Range("J199").Value = 250
$ws.Range("K199").Value = 2000
$ws.Range("L199").Value = 2200
$ws.Range("M199").Value = 2104
$ws.Range("N199").Value = "`$/docena de atados (6 kilos)"
$ws.Range("O199").Value = "Provincia de Quillota"
$ws.Range("P199").Value = 351
$ws.Range("Q199").Value = 6
$ws.Range("R199").Value = "Hortaliza"
